$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.474.48'

$ws.Range('D3').Value = '2.631.96'
$ws.Range('E3').Value = '  -1.46%  '

$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.54%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.33'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.64%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.645'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.93%  '

$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('E9').Value = '  -4.80%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.80'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.23%  '

$ws.Range('E11').Value = '  -2.44%  '

$ws.Range('E12').Value = '  -0.09%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.53'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.73%  '

$ws.Range('E14').Value = '  -6.34%  '

$ws.Range('D15').Value = '3.104.41'
$ws.Range('E15').Value = '  -1.47%  '

$ws.Range('D16').Value = '64.266.56'
$ws.Range('E16').Value = '  -2.08%  '

$ws.Range('D17').Value = '2.634.18'
$ws.Range('E17').Value = '  -1.66%  '

$ws.Range('E18').Value = '  -3.33%  '

$ws.Range('E19').Value = '  -2.14%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.41'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.12%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '346.53'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.46%  '

$ws.Range('E22').Value = '  -0.11%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.72'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.10%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000113'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.26%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.75'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.22%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.35'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.85%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.57'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.40%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '560.01'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.92%  '

$ws.Range('E29').Value = '  -2.11%  '

$ws.Range('E30').Value = '  -0.01%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.91'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.08%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.07'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.98%  '

$ws.Range('E33').Value = '  -3.03%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.60'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.03%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.28'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.58%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.412'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.63%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.05'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.55%  '

$ws.Range('E38').Value = '  -0.04%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.93'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.24%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '154.34'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.63%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.01%  '

$ws.Range('E42').Value = '  +4.78%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '158.76'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.57%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.98'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.78%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0598'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.99%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.84'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.71%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.636'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.61%  '

$ws.Range('E48').Value = '  +3.12%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.16'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.09%  '

$ws.Range('D51').Value = '0.0₆0240'
$ws.Range('E51').Value = '  -5.38%  '

